$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update the "Weights" row (row 7): the Face Coverings weight (AG)
#    is switched off.
# ---------------------------------------------------------------------
$oldTotal = $ws.Cells.Item(7, 35).Value2
$oldAG = $ws.Cells.Item(7, 33).Value2

$ws.Range("AG7").Value = 0

# Recompute the total weight (AI7 = SUM(B7:AH7)); since only AG changed,
# simply remove its old contribution from the previously stored total.
$total = $oldTotal - $oldAG
$ws.Cells.Item(7, 35).Value = $total

# ---------------------------------------------------------------------
# 2) Recompute the LockdownEffectiveness column (AI) for every data
#    row (9 through 221) using the new weights / new total.
#    AI{r} = SUMPRODUCT(B{r}:AH{r}, B7:AH7) / AI7
# ---------------------------------------------------------------------
$weights = @{}
for ($c = 2; $c -le 34; $c++) {
    $weights[$c] = $ws.Cells.Item(7, $c).Value2
}

for ($r = 9; $r -le 221; $r++) {
    $num = 0
    for ($c = 2; $c -le 34; $c++) {
        $num = $num + ($ws.Cells.Item($r, $c).Value2 * $weights[$c])
    }
    $ws.Cells.Item($r, 35).Value = $num / $total
}

# ---------------------------------------------------------------------
# 3) Append 12 new dates (9/30/2020 - 10/11/2020) as new rows 222-233,
#    following the same pattern as the last existing row (221).
# ---------------------------------------------------------------------
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")

$templateValues = @{
    2  = 0;  # B
    3  = 0;  # C
    4  = 1;  # D
    5  = 1;  # E
    6  = 1;  # F
    7  = 0;  # G
    8  = 0;  # H
    9  = 0;  # I
    10 = 0;  # J
    11 = 0;  # K
    12 = 0;  # L
    13 = 0;  # M
    14 = 0;  # N
    15 = 0;  # O
    16 = 0;  # P
    17 = 0;  # Q
    18 = 0;  # R
    19 = 0;  # S
    20 = 0;  # T
    21 = 0;  # U
    22 = 0;  # V
    23 = 0;  # W
    24 = 0;  # X
    25 = 0;  # Y
    26 = 0;  # Z
    27 = 0;  # AA
    28 = 0;  # AB
    29 = 0;  # AC
    30 = 0;  # AD
    31 = 0;  # AE
    32 = 0;  # AF
    33 = 1;  # AG
    34 = 1;  # AH
}

$row = 222
foreach ($d in $newDates) {
    # Force the date-looking text to be stored as a string (not an auto
    # parsed date) by entering it with a leading apostrophe, then copy
    # the formatting (bold / border / centered) from the last templated
    # row so the underlying style stays identical to the existing rows.
    $ws.Cells.Item($row, 1).Value = "'" + $d
    $ws.Range("A221").Copy()
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 1)).PasteSpecial(-4122)  # xlPasteFormats

    $num = 0
    for ($c = 2; $c -le 34; $c++) {
        $v = $templateValues[$c]
        $ws.Cells.Item($row, $c).Value = $v
        $num = $num + ($v * $weights[$c])
    }
    $ws.Cells.Item($row, 35).Value = $num / $total

    $row = $row + 1
}
